# Refresh the cryptocurrency price/volume table to match the
# "Mon May  8 17:18:21 UTC 2023" GitHub Actions data update.
# (Also re-ranks Hedera/TrustWalletToken at rows 37-38 and
# Quant/EnergySwap at rows 47-48, swapping their data.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "1.007") are
# given a leading apostrophe so Excel stores them as literal text -
# exactly like the original inline-string cells - instead of silently
# re-interpreting them as numeric values.

$ws.Range("D2").Value = "28.005.62"
$ws.Range("E2").Value = "  -3.71%  "
$ws.Range("D3").Value = "1.869.52"
$ws.Range("E3").Value = "  -2.95%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'318.69"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.4336"
$ws.Range("E7").Value = "  -5.81%  "
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "'0.07408"
$ws.Range("E9").Value = "  -4.61%  "
$ws.Range("D10").Value = "'0.9304"
$ws.Range("E10").Value = "  -5.22%  "
$ws.Range("D11").Value = "'21.23"
$ws.Range("E11").Value = "  -6.69%  "
$ws.Range("D12").Value = "1.906.31"
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "'6.727"
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "'5.427"
$ws.Range("E14").Value = "  -4.74%  "
$ws.Range("D15").Value = "'0.06871"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "'80.39"
$ws.Range("E17").Value = "  -4.68%  "
$ws.Range("D18").Value = "'0.000009004"
$ws.Range("E18").Value = "  -5.63%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'15.74"
$ws.Range("E20").Value = "  -6.06%  "
$ws.Range("D21").Value = "28.003.92"
$ws.Range("E21").Value = "  -3.77%  "
$ws.Range("D22").Value = "'5.114"
$ws.Range("E22").Value = "  -4.34%  "
$ws.Range("D23").Value = "'11.00"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "2.122.81"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").Value = "'2.048"
$ws.Range("E25").Value = "  -1.38%  "
$ws.Range("D26").Value = "'154.05"
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("D27").Value = "'18.51"
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("D28").Value = "'5.499"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("D29").Value = "'113.17"
$ws.Range("E29").Value = "  -4.14%  "
$ws.Range("D30").Value = "'1.692"
$ws.Range("E30").Value = "  -7.90%  "
$ws.Range("D31").Value = "'0.08978"
$ws.Range("E31").Value = "  -3.98%  "
$ws.Range("D32").Value = "'0.8067"
$ws.Range("E32").Value = "  -6.33%  "
$ws.Range("D33").Value = "'4.767"
$ws.Range("E33").Value = "  -6.90%  "
$ws.Range("D34").Value = "'1.174"
$ws.Range("E34").Value = "  -5.62%  "
$ws.Range("D35").Value = "'2.955"
$ws.Range("E35").Value = "  -1.98%  "
$ws.Range("D36").Value = "'1.005"
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.05498"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.118"
$ws.Range("E38").Value = "  -3.63%  "
$ws.Range("E39").Value = "  -3.85%  "
$ws.Range("D40").Value = "'3.000"
$ws.Range("E40").Value = "  -6.50%  "
$ws.Range("D41").Value = "'0.5237"
$ws.Range("E41").Value = "  -5.07%  "
$ws.Range("D42").Value = "'7.008"
$ws.Range("E42").Value = "  -6.60%  "
$ws.Range("D43").Value = "'0.1685"
$ws.Range("E43").Value = "  -4.12%  "
$ws.Range("D44").Value = "'8.742"
$ws.Range("E44").Value = "  -6.76%  "
$ws.Range("D45").Value = "'0.06703"
$ws.Range("E45").Value = "  -3.08%  "
$ws.Range("D46").Value = "'0.4871"
$ws.Range("E46").Value = "  -6.36%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'106.81"
$ws.Range("E47").Value = "  -3.40%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'10.44"
$ws.Range("E48").Value = "  -7.85%  "
$ws.Range("D49").Value = "'1.003"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'1.669"
$ws.Range("E50").Value = "  -5.74%  "
$ws.Range("D51").Value = "'1.873"
$ws.Range("E51").Value = "  -14.50%  "
